# Scene 50 edit: apply 1.15x ("single-and-a-bit") line spacing to every
# paragraph, and merge the three split runs in the "behaviour's changed"
# line back into a single run (the text itself is unchanged - it had been
# typed/edited across runs).

$d = $word.ActiveDocument

# --- 1. Give every paragraph w:spacing line="276" lineRule="auto" (1.15 lines) ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $para.LineSpacingRule = 5   # wdLineSpaceMultiple
    $para.LineSpacing = 13.8    # 276/20 -> 1.15 * 12pt single spacing
}

# --- 2. Re-merge the runs that make up the "Petra (neutral thinking)" line ---
$apos = [char]0x2019

$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*behaviour*s changed*") {
        $target = $para
        break
    }
}

if ($target -ne $null) {
    $pStart = $target.Range.Start
    $pEnd = $target.Range.End

    # Locate the start of the second run's text ("behaviour's changed ...")
    # inside this paragraph so we can collapse the trailing runs into the
    # first one without disturbing its own run formatting/attributes.
    $scan = $d.Range($pStart, $pEnd)
    $scan.Find.Execute("behaviour" + $apos + "s changed", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

    $splitPoint = $scan.Start
    $tailEnd = $pEnd - 1   # stop before the paragraph mark

    $tail = $d.Range($splitPoint, $tailEnd)
    $tailText = $tail.Text
    $tail.Delete()

    $head = $d.Range($pStart, $splitPoint)
    $head.InsertAfter($tailText)
}

Write-Output "done"
